$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new row of data (row 17) to the tracker table
$ws.Range("A17").Value = "pairwise"
$ws.Range("C17").Value = "RD"
$ws.Range("D17").Value = "all_submitted_tracker_RD_Nov-20-2023.csv"

# "Nov-20-2023" looks like a date to Excel's auto-detection, which would
# turn it into a date serial number. Build it via a text-literal formula in
# a scratch cell, then paste just the resulting value into B17 so it lands
# as plain text (matching the other Date-column entries in this sheet).
$ws.Range("Z1").Formula = '="Nov-20-2023"'
$ws.Range("Z1").Copy()
$ws.Range("B17").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").Clear()
$excel.CutCopyMode = 0

# Move the active selection to D17, matching the saved view state
$ws.Range("D17").Select()
